$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 69
$ws1.Range("F3").Value = 11730
$ws1.Range("F4").Value = 216
$ws1.Range("F5").Value = 339
$ws1.Range("F7").Value = 11675
$ws1.Range("F12").Value = 5779
$ws1.Range("F14").Value = 3518
$ws1.Range("F15").Value = 185

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 571

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 571
$ws4.Range("F3").Value = 69
$ws4.Range("F5").Value = 11730
$ws4.Range("F6").Value = 216
$ws4.Range("F9").Value = 11675
$ws4.Range("F15").Value = 5779
$ws4.Range("F17").Value = 3518
$ws4.Range("F18").Value = 185
